$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''63.632.61'
$ws.Range("E2").Value = '''  -1.70%  '
$ws.Range("D3").Value = '''3.043.04'
$ws.Range("E3").Value = '''  -1.91%  '
$ws.Range("E4").Value = '''  +0.08%  '
$ws.Range("D5").Value = '''557.05'
$ws.Range("E5").Value = '''  -0.40%  '
$ws.Range("D6").Value = '''141.61'
$ws.Range("E6").Value = '''  -1.80%  '
$ws.Range("E7").Value = '''  +0.07%  '
$ws.Range("D8").Value = '''3.040.67'
$ws.Range("E8").Value = '''  -1.89%  '
$ws.Range("D9").Value = '''0.520'
$ws.Range("E9").Value = '''  +4.18%  '
$ws.Range("D10").Value = '''6.34'
$ws.Range("E10").Value = '''  -9.78%  '
$ws.Range("D11").Value = '''0.152'
$ws.Range("E11").Value = '''  +0.37%  '
$ws.Range("D12").Value = '''0.486'
$ws.Range("E12").Value = '''  +5.18%  '
$ws.Range("E13").Value = '''  +0.38%  '
$ws.Range("D14").Value = '''35.51'
$ws.Range("E14").Value = '''  +0.57%  '
$ws.Range("D15").Value = '''3.540.34'
$ws.Range("D16").Value = '''63.682.73'
$ws.Range("E16").Value = '''  -1.50%  '
$ws.Range("D17").Value = '''3.039.00'
$ws.Range("E17").Value = '''  -1.97%  '
$ws.Range("E18").Value = '''  +0.31%  '
$ws.Range("D19").Value = '''6.77'
$ws.Range("E19").Value = '''  -0.98%  '
$ws.Range("D20").Value = '''473.89'
$ws.Range("E20").Value = '''  -2.13%  '
$ws.Range("D21").Value = '''14.03'
$ws.Range("E21").Value = '''  +1.36%  '
$ws.Range("D22").Value = '''14.59'
$ws.Range("E22").Value = '''  +9.64%  '
$ws.Range("D23").Value = '''0.680'
$ws.Range("E23").Value = '''  +0.74%  '
$ws.Range("D24").Value = '''7.50'
$ws.Range("E24").Value = '''  -2.10%  '
$ws.Range("D25").Value = '''82.82'
$ws.Range("E25").Value = '''  +2.28%  '
$ws.Range("D26").Value = '''0.999'
$ws.Range("E26").Value = '''  -0.07%  '
$ws.Range("D27").Value = '''2.78'
$ws.Range("E27").Value = '''  -0.60%  '
$ws.Range("D28").Value = '''8.08'
$ws.Range("E28").Value = '''  -0.21%  '
$ws.Range("D29").Value = '''2.02'
$ws.Range("E29").Value = '''  -2.17%  '
$ws.Range("E30").Value = '''  +0.08%  '
$ws.Range("D31").Value = '''26.14'
$ws.Range("E31").Value = '''  +0.06%  '
$ws.Range("E32").Value = '''  -1.09%  '
$ws.Range("E33").Value = '''  -1.52%  '
$ws.Range("D34").Value = '''5.72'
$ws.Range("E34").Value = '''  -0.01%  '
$ws.Range("D35").Value = '''6.19'
$ws.Range("E35").Value = '''  -0.17%  '
$ws.Range("D36").Value = '''54.62'
$ws.Range("E36").Value = '''  -0.78%  '
$ws.Range("D37").Value = '''0.0407'
$ws.Range("E37").Value = '''  -0.19%  '
$ws.Range("D38").Value = '''439.55'
$ws.Range("E38").Value = '''  -5.49%  '
$ws.Range("D39").Value = '''0.0810'
$ws.Range("E39").Value = '''  -2.33%  '
$ws.Range("D40").Value = '''3.005.43'
$ws.Range("E40").Value = '''  -0.31%  '
$ws.Range("D41").Value = '''2.77'
$ws.Range("E41").Value = '''  +1.47%  '
$ws.Range("E42").Value = '''  -0.01%  '
$ws.Range("D43").Value = '''8.25'
$ws.Range("E43").Value = '''  -0.34%  '
$ws.Range("D44").Value = '''0.268'
$ws.Range("E44").Value = '''  +2.66%  '
$ws.Range("D45").Value = '''27.72'
$ws.Range("E45").Value = '''  -3.27%  '
$ws.Range("D46").Value = '''2.23'
$ws.Range("E46").Value = '''  +6.17%  '
$ws.Range("E47").Value = '''  -0.02%  '
$ws.Range("D48").Value = '''0.114'
$ws.Range("E48").Value = '''  +1.05%  '
$ws.Range("D49").Value = '''118.19'
$ws.Range("E49").Value = '''  -0.61%  '
$ws.Range("D50").Value = '''0.0₃0511'
$ws.Range("E50").Value = '''  -0.85%  '
$ws.Range("D51").Value = '''2.07'
$ws.Range("E51").Value = '''  -0.01%  '
